$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that needs to move
# from 45203 (2023-10-04) to 45205 (2023-10-06) for every data row (2-171).
$ws.Range("C2:C171").Value = 45205
